$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a last-modified date serial for each data row.
# Bump it from 45205 (2023-10-06) to 45206 (2023-10-07) for every data row (2-46).
for ($row = 2; $row -le 46; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value = 45206
    }
}
